# Daily attendance processing - 2026-01-18 18:41:37
# Normalizes the "Recorded By" (column G) values: for entries that list
# exactly two recorders separated by ", ", swap their order
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"),
# except where one of the recorders is the backup/backdoor account,
# whose order is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    $parts = $text -split ", "

    if ($parts.Count -eq 2 -and $text -notmatch "backup@backdoor.com") {
        $swapped = $parts[1] + ", " + $parts[0]
        $cell.Value = $swapped
    }
}
